$wb = $excel.ActiveWorkbook

# --- 1. Add a "State" column to the hotel_info sheet, between Hotel_Name and City ---
$hotel = $wb.Worksheets.Item("hotel_info")
$hotel.Columns.Item(3).Insert()
$hotel.Cells.Item(1, 3).Value = "State"
$hotel.Cells.Item(2, 3).Value = "Louisiana"

# --- 2. Reorder the sheet tabs so review_info comes before hotel_info ---
$review = $wb.Worksheets.Item("review_info")
$review.Move($hotel)

$wb.Save()
